$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new entry "005" ---
$ws.Rows(6).RowHeight = 221.25

$ws.Range("A6").Value = "'005"
$ws.Range("A6").VerticalAlignment = -4160

$texto005 = "+ palancas TRUE: corregir, nuevasvars, lag1, delta1, lag2, delta2, ratiomax3, deflactar, canaritosimportancia`n+ variables drift: `"mpasivos_margen`", `"mactivos_margen`", `"mcuentas_saldo`",`n                              `"mcajeros_propios_descuentos`", `"mtarjeta_visa_descuentos`",`n                              `"mforex_sell`", `"mtransferencias_emitidas`", `n                              `"Master_mfinanciacion_limite`",`"Master_mconsumospesos`",`n                              `"Master_fultimo_cierre`", `"Master_madelantodolares`",`"Master_mpagado`",`n                              `"Master_mpagominimo`", `"Master_mconsumototal`", `n                              `"Visa_mfinanciacion_limite`",`n                              `"Visa_msaldototal`", `"Visa_msaldopesos`", `"Visa_msaldodolares`",`n                              `"Visa_mconsumospesos`", `"Visa_fultimo_cierre`", `"Visa_mconsumototal`",`n                              `"Visa_mpagominimo`"`n+ el resto es todo igual al 004"

$ws.Range("B6").Value = "'" + $texto005
$ws.Range("B6").WrapText = $true
$ws.Range("B6").HorizontalAlignment = -4131
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("B6:I6").Merge()

# --- Row 7: new entry "006" ---
$ws.Range("A7").Value = "'006"
$ws.Range("A7").VerticalAlignment = -4160

$ws.Range("B7").Value = "igual al 005 pero sin DEFLACTAR!"
$ws.Range("B7").VerticalAlignment = -4160

# --- Scroll / selection, mirroring where the author's cursor ended up ---
[void]$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 6

Write-Output "done"
